$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column C
$ws.Range("C1").Value = "Z"

# Copy column A (rows 2-7) into column C
for ($r = 2; $r -le 7; $r++) {
    $aText = $ws.Cells.Item($r, 1).Text
    $ws.Cells.Item($r, 3).Value = $aText
}

# Update the active selection to mirror the saved view state (C8)
$ws.Range("C8").Select()
